$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 21-67: existing rows being re-keyed / re-valued for new shapes categories (ellipse, blob)
# and re-ordered underlying data. Column C score weight is updated where it changed.

$ws.Cells.Item(21,1).Value = "ellipse"
$ws.Cells.Item(21,2).Value = "Ellipse"

$ws.Cells.Item(22,1).Value = "blob"
$ws.Cells.Item(22,2).Value = "Blob"
$ws.Cells.Item(22,3).ClearContents()

$ws.Cells.Item(23,1).Value = "polygon"
$ws.Cells.Item(23,2).Value = "Polygon"
$ws.Cells.Item(23,3).ClearContents()

$ws.Cells.Item(24,1).Value = "triangle"
$ws.Cells.Item(24,2).Value = "Triangle"
$ws.Cells.Item(24,3).Value = 1.5

$ws.Cells.Item(25,1).Value = "quadrilateral"
$ws.Cells.Item(25,2).Value = "Quadrilateral"
$ws.Cells.Item(25,3).Value = 1.5

$ws.Cells.Item(26,1).Value = "pentagon"
$ws.Cells.Item(26,2).Value = "Pentagon"
$ws.Cells.Item(26,3).Value = 1.5

$ws.Cells.Item(27,1).Value = "hexagon"
$ws.Cells.Item(27,2).Value = "Hexagon"
$ws.Cells.Item(27,3).Value = 1.5

$ws.Cells.Item(28,1).Value = "octagon"
$ws.Cells.Item(28,2).Value = "Octagon"
$ws.Cells.Item(28,3).Value = 1.5

$ws.Cells.Item(29,1).Value = "decagon"
$ws.Cells.Item(29,2).Value = "Decagon"
$ws.Cells.Item(29,3).Value = 1.5

$ws.Cells.Item(30,1).Value = "triangle_right"
$ws.Cells.Item(30,2).Value = "Right Triangle"
$ws.Cells.Item(30,3).Value = 2

$ws.Cells.Item(31,1).Value = "triangle_equilateral"
$ws.Cells.Item(31,2).Value = "Equilateral Triangle"
$ws.Cells.Item(31,3).Value = 2

$ws.Cells.Item(32,1).Value = "triangle_isosceles"
$ws.Cells.Item(32,2).Value = "Isosceles Triangle"
$ws.Cells.Item(32,3).Value = 2

$ws.Cells.Item(33,1).Value = "triangle_scalene"
$ws.Cells.Item(33,2).Value = "Scalene Triangle"
$ws.Cells.Item(33,3).Value = 2

$ws.Cells.Item(34,1).Value = "triangle_acute"
$ws.Cells.Item(34,2).Value = "Acute Triangle"
$ws.Cells.Item(34,3).Value = 2

$ws.Cells.Item(35,1).Value = "triangle_obtuse"
$ws.Cells.Item(35,2).Value = "Obtuse Triangle"
$ws.Cells.Item(35,3).Value = 2

$ws.Cells.Item(36,1).Value = "quad_rectangle"
$ws.Cells.Item(36,2).Value = "Rectangle"
$ws.Cells.Item(36,3).Value = 2

$ws.Cells.Item(37,1).Value = "quad_square"
$ws.Cells.Item(37,2).Value = "Square"
$ws.Cells.Item(37,3).Value = 2

$ws.Cells.Item(38,1).Value = "quad_rhombus"
$ws.Cells.Item(38,2).Value = "Rhombus"
$ws.Cells.Item(38,3).Value = 2

$ws.Cells.Item(39,1).Value = "quad_parallelogram"
$ws.Cells.Item(39,2).Value = "Parallelogram"
$ws.Cells.Item(39,3).Value = 3

$ws.Cells.Item(40,1).Value = "quad_trapezoid"
$ws.Cells.Item(40,2).Value = "Trapezoid"
$ws.Cells.Item(40,3).Value = 2.5

$ws.Cells.Item(41,1).Value = "quad_trapezium"
$ws.Cells.Item(41,2).Value = "Trapezium"
$ws.Cells.Item(41,3).Value = 2.5

$ws.Cells.Item(42,1).Value = "quad_kite"
$ws.Cells.Item(42,2).Value = "Kite"
$ws.Cells.Item(42,3).Value = 1

$ws.Cells.Item(43,1).Value = "prop_side_3"
$ws.Cells.Item(43,2).Value = "3 sides."
$ws.Cells.Item(43,3).Value = 2

$ws.Cells.Item(44,1).Value = "prop_side_4"
$ws.Cells.Item(44,2).Value = "4 sides."
$ws.Cells.Item(44,3).Value = 2

$ws.Cells.Item(45,1).Value = "prop_side_5"
$ws.Cells.Item(45,2).Value = "5 sides."
$ws.Cells.Item(45,3).Value = 2

$ws.Cells.Item(46,1).Value = "prop_side_6"
$ws.Cells.Item(46,2).Value = "6 sides."
$ws.Cells.Item(46,3).Value = 2

$ws.Cells.Item(47,1).Value = "prop_side_8"
$ws.Cells.Item(47,2).Value = "8 sides."
$ws.Cells.Item(47,3).Value = 2

$ws.Cells.Item(48,1).Value = "prop_side_10"
$ws.Cells.Item(48,2).Value = "10 sides."
$ws.Cells.Item(48,3).Value = 2

$ws.Cells.Item(49,1).Value = "prop_90_degree"
$ws.Cells.Item(49,2).Value = "Has a 90° angle."
$ws.Cells.Item(49,3).Value = 3

$ws.Cells.Item(50,1).Value = "prop_sides_equal_all"
$ws.Cells.Item(50,2).Value = "All sides are equal."
$ws.Cells.Item(50,3).Value = 3

$ws.Cells.Item(51,1).Value = "prop_60_degree_all"
$ws.Cells.Item(51,2).Value = "All angles equal 60°."
$ws.Cells.Item(51,3).Value = 3

$ws.Cells.Item(52,1).Value = "prop_sides_equal_two"
$ws.Cells.Item(52,2).Value = "Two equal sides."
$ws.Cells.Item(52,3).Value = 3

$ws.Cells.Item(53,1).Value = "prop_angles_equal_two"
$ws.Cells.Item(53,2).Value = "Two equal angles."
$ws.Cells.Item(53,3).Value = 3

$ws.Cells.Item(54,1).Value = "prop_sides_no_equal"
$ws.Cells.Item(54,2).Value = "No sides are equal."
$ws.Cells.Item(54,3).Value = 3

$ws.Cells.Item(55,1).Value = "prop_angles_less_90_all"
$ws.Cells.Item(55,2).Value = "All angles less than 90°."
$ws.Cells.Item(55,3).Value = 3

$ws.Cells.Item(56,1).Value = "prop_angle_greater_90"
$ws.Cells.Item(56,2).Value = "Has an angle greater than 90°."
$ws.Cells.Item(56,3).Value = 3

$ws.Cells.Item(57,1).Value = "prop_90_degree_all"
$ws.Cells.Item(57,2).Value = "All angles equal 90°."
$ws.Cells.Item(57,3).Value = 3

$ws.Cells.Item(58,1).Value = "prop_opposite_sides_parallel"
$ws.Cells.Item(58,2).Value = "Opposite sides are parallel."
$ws.Cells.Item(58,3).Value = 3

$ws.Cells.Item(59,1).Value = "prop_opposite_sides_equal"
$ws.Cells.Item(59,2).Value = "Opposite sides are equal."
$ws.Cells.Item(59,3).Value = 3

$ws.Cells.Item(60,1).Value = "prop_opposite_angles_equal"
$ws.Cells.Item(60,2).Value = "Opposite angles are equal."
$ws.Cells.Item(60,3).Value = 3

$ws.Cells.Item(61,1).Value = "prop_opposite_pair_parallel"
$ws.Cells.Item(61,2).Value = "Two sides are parallel."
$ws.Cells.Item(61,3).Value = 3

$ws.Cells.Item(62,1).Value = "prop_sides_no_parallel"
$ws.Cells.Item(62,2).Value = "No sides are parallel."
$ws.Cells.Item(62,3).Value = 3

$ws.Cells.Item(63,1).Value = "prop_sides_pair_equal"
$ws.Cells.Item(63,2).Value = "Two pairs of sides are equal."
$ws.Cells.Item(63,3).Value = 3

$ws.Cells.Item(64,1).Value = "shape_categories"
$ws.Cells.Item(64,2).Value = "Shape Categories"

$ws.Cells.Item(65,1).Value = "shapes"
$ws.Cells.Item(65,2).Value = "Shapes"
$ws.Cells.Item(65,2).WrapText = $true

$ws.Cells.Item(66,1).Value = "proceed_instruct"
$ws.Cells.Item(66,2).Value = "Press this button to proceed."
$ws.Cells.Item(66,3).Value = 3
$ws.Cells.Item(66,2).WrapText = $true

$ws.Cells.Item(67,1).Value = "victory"
$ws.Cells.Item(67,2).Value = "VICTORY"
$ws.Cells.Item(67,2).WrapText = $true

$ws.Cells.Item(68,1).Value = "ellipse_desc"
$ws.Cells.Item(68,2).Value = "· Has a center point.\n\n· Dimensions defined by two axis."

$ws.Cells.Item(69,1).Value = "blob_desc"
$ws.Cells.Item(69,2).Value = "· Made up of curves defined by equations."

$ws.Cells.Item(70,1).Value = "poly_desc"
$ws.Cells.Item(70,2).Value = "· Has three or more points plotted on a plane.\n\n· Has a number of sides connected by points."

$ws.Cells.Item(71,1).Value = "side"
$ws.Cells.Item(71,2).Value = "Side (length)"

$ws.Cells.Item(72,1).Value = "point"
$ws.Cells.Item(72,2).Value = "Point (x, y)"

$ws.Cells.Item(73,1).Value = "instruct_collect"
$ws.Cells.Item(73,2).Value = "Press and hold the shape to collect it."

$ws.Cells.Item(74,1).Value = "instruct_category_drag"
$ws.Cells.Item(74,2).Value = "Press and drag the item towards the category, and release."

$ws.Cells.Item(75,1).Value = "intro_0_0"
$ws.Cells.Item(75,2).Value = "Welcome to Kero Shapes!"

$ws.Cells.Item(76,1).Value = "intro_0_1"
$ws.Cells.Item(76,2).Value = "In this game, you will be helping these fine capable frogs to categorize shapes that are scattered about in the land."

$ws.Cells.Item(77,1).Value = "intro_0_2"
$ws.Cells.Item(77,2).Value = "When it comes to categorizing shapes, one must figure out the common attributes they share."

$ws.Cells.Item(78,1).Value = "intro_1_0"
$ws.Cells.Item(78,2).Value = "Here are some primary categories and their attributes that differentiate them from each other."

$ws.Cells.Item(79,1).Value = "intro_1_1"
$ws.Cells.Item(79,2).Value = "For now, we will be focusing on identifying polygons."

$ws.Cells.Item(80,1).Value = "intro_2_0"
$ws.Cells.Item(80,2).Value = "Remember that a polygon is formed by three or more straight lines that are connected in a loop."

$ws.Cells.Item(81,1).Value = "intro_3_0"
$ws.Cells.Item(81,2).Value = "Now let's go forth, and help these lads survey the land!"

$ws.Cells.Item(82,1).Value = "level_intro_0_0"
$ws.Cells.Item(82,2).Value = "The following types of shapes will come up. Be sure to remember them!"

$ws.Cells.Item(83,1).Value = "level_0_continue"
$ws.Cells.Item(83,2).Value = "Go ahead and collect the rest of the shapes! Remember to match them with the correct category."

# Update the active selection / scroll position to match the end of the newly-entered data.
$ws.Range("B83").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 64
$excel.ActiveWindow.ScrollColumn = 1
